# Generate Report for Handback
# - Flip Status from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview rollup sheet (columns zh-cn/de-de) and on each language
#   sheet's Status column.
# - Record the handback results on each language sheet: Latest Target File
#   (hyperlinked to the source file), Latest Handback File (the generated
#   xliff), and Latest Handback DateTime (the timestamp of the handback).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ffd4a5f8907ba81940a99670674925bec24cdad/e2e/"

# ---- Overview sheet: flip the rolled-up status for both files/languages ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

function Set-Handback {
    param(
        [string]$SheetName,
        [string]$XliffName,
        [string]$HandbackTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (I): hyperlink to a.md, same as column A's link.
    $ws.Hyperlinks.Add($ws.Range("I2"), ($ghBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), ($ghBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J): the generated xliff file name.
    $ws.Range("J2").Value = $XliffName
    $ws.Range("J3").Value = $XliffName

    # Latest Handback DateTime (K): when the handback was generated.
    $ws.Range("K2").Value = $HandbackTime
    $ws.Range("K3").Value = $HandbackTime
}

Set-Handback "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-15 18:33:37"
Set-Handback "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-15 18:33:45"
